# New crime data collected: weekly CompStat figures updated (Volume 32 No. 48, week 11/24/2025-11/30/2025)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text cells: Volume/Number and reporting week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Simple value-only updates (style unchanged) ---
$ws.Range("N14").Value = -94.666666666666
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = 27.777777777777
$ws.Range("N15").Value = -65.671641791044
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -31.578947368421
$ws.Range("I16").Value = 172
$ws.Range("J16").Value = 184
$ws.Range("K16").Value = -6.521739130434
$ws.Range("L16").Value = -16.504854368932
$ws.Range("M16").Value = -32.015810276679
$ws.Range("N16").Value = -91.404297851074
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -38.461538461538
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -14.634146341463
$ws.Range("I17").Value = 425
$ws.Range("J17").Value = 370
$ws.Range("K17").Value = 14.864864864864
$ws.Range("L17").Value = 8.974358974358
$ws.Range("M17").Value = 50.709219858156
$ws.Range("N17").Value = -61.677186654643
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("H18").Value = -41.666666666666
$ws.Range("J18").Value = 139
$ws.Range("K18").Value = 3.597122302158
$ws.Range("L18").Value = -17.241379310344
$ws.Range("M18").Value = -50.684931506849
$ws.Range("N18").Value = -90.123456790123
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 25.806451612903
$ws.Range("I19").Value = 415
$ws.Range("J19").Value = 339
$ws.Range("K19").Value = 22.418879056047
$ws.Range("L19").Value = 10.962566844919
$ws.Range("M19").Value = 49.280575539568
$ws.Range("N19").Value = -23.572744014733
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 300
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 133.333333333333
$ws.Range("I20").Value = 113
$ws.Range("J20").Value = 99
$ws.Range("K20").Value = 14.141414141414
$ws.Range("L20").Value = -26.623376623376
$ws.Range("M20").Value = -5.042016806722
$ws.Range("N20").Value = -83.857142857142
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -14.814814814814
$ws.Range("F21").Value = 109
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -0.909090909090
$ws.Range("I21").Value = 1296
$ws.Range("J21").Value = 1164
$ws.Range("K21").Value = 11.340206185567
$ws.Range("L21").Value = -2.040816326530
$ws.Range("M21").Value = 2.775574940523
$ws.Range("N21").Value = -78.229464135729
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -83.333333333333
$ws.Range("I22").Value = 22
$ws.Range("J22").Value = 25
$ws.Range("K22").Value = -12
$ws.Range("L22").Value = -18.518518518518
$ws.Range("M22").Value = -35.294117647058
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -22.222222222222
$ws.Range("I23").Value = 105
$ws.Range("J23").Value = 109
$ws.Range("K23").Value = -3.669724770642
$ws.Range("L23").Value = -1.869158878504
$ws.Range("M23").Value = 59.090909090909
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -12.5
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 56
$ws.Range("H24").Value = 7.142857142857
$ws.Range("I24").Value = 773
$ws.Range("J24").Value = 688
$ws.Range("K24").Value = 12.354651162790
$ws.Range("L24").Value = 9.335219236209
$ws.Range("M24").Value = -1.277139208173
$ws.Range("C25").Value = 3
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 57.142857142857
$ws.Range("I25").Value = 94
$ws.Range("K25").Value = -16.071428571428
$ws.Range("L25").Value = -1.052631578947
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = 4.761904761904
$ws.Range("I26").Value = 411
$ws.Range("J26").Value = 467
$ws.Range("K26").Value = -11.991434689507
$ws.Range("L26").Value = -20.502901353965
$ws.Range("M26").Value = -44.906166219839
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 50
$ws.Range("L27").Value = 3.225806451612
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 79
$ws.Range("J28").Value = 51
$ws.Range("K28").Value = 54.901960784313
$ws.Range("L28").Value = 17.910447761194
$ws.Range("L31").Value = -40

# --- Cells changing to numeric with a different style: paste format from a same-style donor cell, then set value ---
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = 1
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = 0
$ws.Range("I14").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = 1

# --- Cells changing to shared/placeholder text ("0" or "***.*"): set text value first (apostrophe forces text), then paste format from a same-style donor cell ---
$ws.Range("G14").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("D25").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4122) | Out-Null
$ws.Range("H29").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Value = "'0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Range("H30").Value = "'***.*"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null

# --- Column width tweaks for columns I (9) and J (10): narrow to match columns C/D/F/G ---
$refWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(9).ColumnWidth = $refWidth
$ws.Columns.Item(10).ColumnWidth = $refWidth

$excel.CutCopyMode = $false